$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue 2 4 "69.319.04"
Set-TextValue 2 5 "  -0.38%  "

# Row 3
Set-TextValue 3 4 "3.433.84"
Set-TextValue 3 5 "  +1.17%  "

# Row 4
Set-TextValue 4 4 "1.00"
Set-TextValue 4 5 "  +0.13%  "

# Row 5
Set-TextValue 5 4 "580.18"
Set-TextValue 5 5 "  -1.49%  "

# Row 6
Set-TextValue 6 4 "177.21"
Set-TextValue 6 5 "  -2.37%  "

# Row 7
Set-TextValue 7 4 "3.428.94"
Set-TextValue 7 5 "  +1.10%  "

# Row 8
Set-TextValue 8 5 "  +0.07%  "

# Row 9
Set-TextValue 9 4 "0.592"
Set-TextValue 9 5 "  -0.68%  "

# Row 10
Set-TextValue 10 4 "0.197"
Set-TextValue 10 5 "  +0.24%  "

# Row 11
Set-TextValue 11 4 "0.586"
Set-TextValue 11 5 "  -0.56%  "

# Row 12
Set-TextValue 12 4 "48.68"
Set-TextValue 12 5 "  -0.42%  "

# Row 13
Set-TextValue 13 5 "  -0.85%  "

# Row 14
Set-TextValue 14 4 "702.34"
Set-TextValue 14 5 "  +1.36%  "

# Row 15
Set-TextValue 15 4 "3.966.64"
Set-TextValue 15 5 "  +0.51%  "

# Row 16
Set-TextValue 16 5 "  +0.61%  "

# Row 17
Set-TextValue 17 4 "69.424.99"
Set-TextValue 17 5 "  -0.22%  "

# Row 18
Set-TextValue 18 4 "3.425.81"
Set-TextValue 18 5 "  +0.92%  "

# Row 19
Set-TextValue 19 5 "  +0.79%  "

# Row 20
Set-TextValue 20 5 "  +0.15%  "

# Row 21
Set-TextValue 21 4 "11.44"
Set-TextValue 21 5 "  +0.14%  "

# Row 22
Set-TextValue 22 5 "  -0.47%  "

# Row 23
Set-TextValue 23 4 "5.39"
Set-TextValue 23 5 "  -0.47%  "

# Row 24
Set-TextValue 24 5 "  -0.96%  "

# Row 25
Set-TextValue 25 4 "101.27"
Set-TextValue 25 5 "  -3.16%  "

# Row 26
Set-TextValue 26 5 "  -2.11%  "

# Row 27
Set-TextValue 27 4 "2.67"
Set-TextValue 27 5 "  -2.03%  "

# Row 28
Set-TextValue 28 4 "9.62"
Set-TextValue 28 5 "  -0.53%  "

# Row 29
Set-TextValue 29 4 "33.66"
Set-TextValue 29 5 "  -2.23%  "

# Row 30
Set-TextValue 30 4 "8.76"
Set-TextValue 30 5 "  +0.58%  "

# Row 31
Set-TextValue 31 5 "  -1.22%  "

# Row 32
Set-TextValue 32 2 "Bittensor"
Set-TextValue 32 3 "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue 32 4 "569.50"
Set-TextValue 32 5 "  +1.88%  "

# Row 33
Set-TextValue 33 2 "dogwifhat"
Set-TextValue 33 3 "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue 33 4 "3.73"
Set-TextValue 33 5 "  +1.34%  "

# Row 34
Set-TextValue 34 4 "11.03"
Set-TextValue 34 5 "  -1.45%  "

# Row 35
Set-TextValue 35 5 "  -1.66%  "

# Row 36
Set-TextValue 36 4 "58.11"
Set-TextValue 36 5 "  -0.69%  "

# Row 37
Set-TextValue 37 4 "0.999"
Set-TextValue 37 5 "  -0.04%  "

# Row 38
Set-TextValue 38 4 "3.605.84"
Set-TextValue 38 5 "  -3.40%  "

# Row 39
Set-TextValue 39 5 "  -1.99%  "

# Row 40
Set-TextValue 40 4 "34.93"
Set-TextValue 40 5 "  -0.27%  "

# Row 41
Set-TextValue 41 4 "0.0₃0733"
Set-TextValue 41 5 "  +3.02%  "

# Row 42
Set-TextValue 42 4 "3.28"
Set-TextValue 42 5 "  +1.18%  "

# Row 43
Set-TextValue 43 5 "  +0.07%  "

# Row 44
Set-TextValue 44 2 "TheGraph"
Set-TextValue 44 3 "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue 44 4 "0.333"
Set-TextValue 44 5 "  -2.27%  "

# Row 45
Set-TextValue 45 2 "VeChain"
Set-TextValue 45 3 "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue 45 4 "0.0420"
Set-TextValue 45 5 "  +0.10%  "

# Row 46
Set-TextValue 46 2 "Mantle"
Set-TextValue 46 3 "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue 46 4 "1.47"
Set-TextValue 46 5 "  +3.76%  "

# Row 47
Set-TextValue 47 2 "ThetaToken"
Set-TextValue 47 3 "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue 47 4 "2.66"
Set-TextValue 47 5 "  -0.10%  "

# Row 48
Set-TextValue 48 2 "Stellar"
Set-TextValue 48 3 "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue 48 4 "0.128"
Set-TextValue 48 5 "  -1.30%  "

# Row 49
Set-TextValue 49 2 "FirstDigitalUSD"
Set-TextValue 49 3 "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue 49 4 "0.999"
Set-TextValue 49 5 "  -0.11%  "

# Row 50
Set-TextValue 50 2 "Monero"
Set-TextValue 50 3 "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue 50 4 "131.41"
Set-TextValue 50 5 "  -0.74%  "

# Row 51
Set-TextValue 51 2 "CoreDAO"
Set-TextValue 51 3 "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
Set-TextValue 51 4 "2.65"
Set-TextValue 51 5 "  +0.87%  "
